$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number 51 -> 52, and week dates ---
$ws.Range("A8").Value = "Volume 31   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/23/2024  Through  12/29/2024"

# --- Convert some numeric cells into text placeholder cells ("0" / "***.*") ---
# Paste VALUES ONLY first (this carries over the shared-string text & type while
# leaving the destination cell style untouched), then paste FORMATS ONLY from the
# same kind of template cell so the style index matches (reuses style 13, no new
# style gets allocated).
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -83.333333333333
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = -88.888888888888
$ws.Range("J16").Value = 163
$ws.Range("K16").Value = 6.134969325153
$ws.Range("L16").Value = -9.424083769633
$ws.Range("M16").Value = -30.522088353413
$ws.Range("N16").Value = -83.381364073006
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 316
$ws.Range("J17").Value = 284
$ws.Range("K17").Value = 11.267605633802
$ws.Range("L17").Value = 6.397306397306
$ws.Range("M17").Value = 68.085106382978
$ws.Range("N17").Value = -38.640776699029
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -6.666666666666
$ws.Range("I18").Value = 169
$ws.Range("J18").Value = 135
$ws.Range("K18").Value = 25.185185185185
$ws.Range("L18").Value = -32.669322709163
$ws.Range("M18").Value = -30.165289256198
$ws.Range("N18").Value = -85.405872193437
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -35.483870967741
$ws.Range("I19").Value = 450
$ws.Range("J19").Value = 557
$ws.Range("K19").Value = -19.210053859964
$ws.Range("L19").Value = -43.109987357775
$ws.Range("M19").Value = 20.967741935483
$ws.Range("N19").Value = -3.01724137931
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -75
$ws.Range("J20").Value = 182
$ws.Range("K20").Value = -20.87912087912
$ws.Range("L20").Value = 3.597122302158
$ws.Range("M20").Value = 21.008403361344
$ws.Range("N20").Value = -85.123966942148
$ws.Range("C21").Value = 7
$ws.Range("E21").Value = -72
$ws.Range("F21").Value = 59
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = -41.584158415841
$ws.Range("I21").Value = 1282
$ws.Range("J21").Value = 1341
$ws.Range("K21").Value = -4.399701715137
$ws.Range("L21").Value = -24.543849323131
$ws.Range("M21").Value = 7.101086048454
$ws.Range("N21").Value = -69.490718705378
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 35
$ws.Range("K22").Value = -31.428571428571
$ws.Range("L22").Value = -31.428571428571
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -45.454545454545
$ws.Range("F24").Value = 61
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = -32.967032967033
$ws.Range("I24").Value = 944
$ws.Range("J24").Value = 1227
$ws.Range("K24").Value = -23.064384678076
$ws.Range("L24").Value = -30.125832716506
$ws.Range("M24").Value = 20.254777070063
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = -58
$ws.Range("I25").Value = 385
$ws.Range("J25").Value = 686
$ws.Range("K25").Value = -43.877551020408
$ws.Range("L25").Value = -48.04318488529
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -45.454545454545
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = -5.128205128205
$ws.Range("I26").Value = 570
$ws.Range("J26").Value = 519
$ws.Range("K26").Value = 9.826589595375
$ws.Range("L26").Value = 28.668171557562
$ws.Range("M26").Value = -4.040404040404
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("L28").Value = -10.76923076923
$ws.Range("N29").Value = -94.117647058823
$ws.Range("N30").Value = -93.75
